$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 1.54

# Row 4
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 9.2
$ws.Range("P4").Value = 4.19
$ws.Range("U4").Value = 1.64
$ws.Range("V4").Value = 2.13

# Row 5
$ws.Range("I5").Value = 2.62

# Row 6
$ws.Range("G6").Value = 1.47
$ws.Range("H6").Value = 3.9
$ws.Range("I6").Value = 6.4
$ws.Range("J6").Value = 2.02
$ws.Range("K6").Value = 2.18
$ws.Range("L6").Value = 5.9
$ws.Range("N6").Value = 11
$ws.Range("O6").Value = 1.24
$ws.Range("P6").Value = 3.3
$ws.Range("Q6").Value = 1.72
$ws.Range("R6").Value = 1.88
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.78
$ws.Range("W6").Value = 6.6
$ws.Range("X6").Value = 6.8
$ws.Range("Z6").Value = 10.25
$ws.Range("AA6").Value = 12
$ws.Range("AC6").Value = 11
$ws.Range("AD6").Value = 7.8
$ws.Range("AE6").Value = 17
$ws.Range("AF6").Value = 80
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 17
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 19.5
$ws.Range("AK6").Value = 150
$ws.Range("AL6").Value = 65
$ws.Range("AM6").Value = 60
$ws.Range("AN6").Value = 3.25
$ws.Range("AO6").Value = 7.1
$ws.Range("AP6").Value = 17.5
$ws.Range("AQ6").Value = 22
$ws.Range("AR6").Value = 55
$ws.Range("AT6").Value = 2.6
$ws.Range("AU6").Value = 7.8
$ws.Range("AV6").Value = 75
$ws.Range("AX6").Value = 7.4
$ws.Range("AY6").Value = 35
$ws.Range("AZ6").Value = 37
$ws.Range("BA6").Value = 250
$ws.Range("BB6").Value = 250
$ws.Range("BC6").Value = 500

# Row 7
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 2.18
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 5.1
$ws.Range("O7").Value = 1.34
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.65
$ws.Range("U7").Value = 1.98
$ws.Range("V7").Value = 1.65
$ws.Range("W7").Value = 5.8
$ws.Range("Y7").Value = 8.5
$ws.Range("Z7").Value = 12
$ws.Range("AA7").Value = 14.5
$ws.Range("AC7").Value = 8.5
$ws.Range("AD7").Value = 6.9
$ws.Range("AE7").Value = 18.5
$ws.Range("AF7").Value = 110
$ws.Range("AH7").Value = 11.75
$ws.Range("AI7").Value = 28
$ws.Range("AJ7").Value = 16.5
$ws.Range("AK7").Value = 90
$ws.Range("AL7").Value = 55
$ws.Range("AM7").Value = 65
$ws.Range("AO7").Value = 7.8
$ws.Range("AQ7").Value = 26
$ws.Range("AT7").Value = 2.57
$ws.Range("AU7").Value = 7.7
$ws.Range("AV7").Value = 75
$ws.Range("AX7").Value = 6.5
$ws.Range("AY7").Value = 28
$ws.Range("BC7").Value = 450

# Row 9
$ws.Range("G9").Value = 1.9
$ws.Range("H9").Value = 3.4
$ws.Range("J9").Value = 2.6
$ws.Range("K9").Value = 2.05
$ws.Range("L9").Value = 4.75
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 2.75
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.58
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.72
$ws.Range("W9").Value = 6
$ws.Range("Y9").Value = 9
$ws.Range("AC9").Value = 8
$ws.Range("AG9").Value = 401
$ws.Range("AY9").Value = 23

# Row 10
$ws.Range("Q10").Value = 1.67
$ws.Range("U10").Value = 1.63
